$wb = $excel.ActiveWorkbook

# Sheets 1-3 and 5: header labels B1:E1 get an "Ano " prefix (e.g. "2015" -> "Ano 2015")
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($sheetName in $anoSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($col in @("B", "C", "D", "E")) {
        $cell = $ws.Range($col + "1")
        $orig = $cell.Value2
        $cell.Value = "Ano " + $orig
    }
}

# Sheet 4: header labels B1:E1 get an "Intervalo " prefix
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
foreach ($col in @("B", "C", "D", "E")) {
    $cell = $ws4.Range($col + "1")
    $orig = $cell.Value2
    $cell.Value = "Intervalo " + $orig
}

# Sheet 6: only B1 gets the "Ano " prefix
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$cell6 = $ws6.Range("B1")
$orig6 = $cell6.Value2
$cell6.Value = "Ano " + $orig6
